$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$pairs = @(
    @("überhaupt", "en général"),
    @("überhaupt nicht", "pas du tout"),
    @("verhindern", "empêcher"),
    @("derzeit", "actuellement"),
    @("ausschliesslich", "exclusif"),
    @("das Gefängnis", "la prison"),
    @("der Knast", "la taule"),
    @("schaffen", "réussir"),
    @("in alle Richtungen spritzen", "gicler dans tous les sens"),
    @("das Waschbecken", "le lavabo"),
    @("ausprobieren", "essayer (mettre à l'épreuve)"),
    @("verkalkt", "entartré")
)

$startRow = 294
for ($i = 0; $i -lt $pairs.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $pairs[$i][0]
    $ws.Cells.Item($row, 2).Value = $pairs[$i][1]
}

$nextRow = $startRow + $pairs.Length
$ws.Range("A$nextRow").Select()
$excel.ActiveWindow.ScrollRow = 286
